$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 value from 10 to 15
$ws.Range("B1").Value = 15

# Add new row 7: Tang, 5, 5 (B7 and C7 stored as text, like B4/C4 etc.)
$ws.Range("A7").Value = "Tang"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "5"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "5"
